# SFC Country Group: give China's SARs (Hong Kong, Macau, Taiwan) and
# Mainland China their own distinct "SFC Country Group" text instead of
# re-using the plain BBG country name, so lookups/grouping by this column
# start to work as separate buckets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 -> Hong Kong
$ws.Range("C20").Value = "China - Hong Kong "

# Row 15 -> China (Mainland)
$ws.Range("C15").Value = "China - Mainland"

# Row 23 -> Macau
$ws.Range("C23").Value = "China - Macau"

# Row 30 -> Taiwan
$ws.Range("C30").Value = "China - Taiwan"

# Restore the active cell the sheet was saved with (the view was also
# scrolled down to row 7, but the headless host doesn't expose/persist
# plain scroll position outside of frozen panes).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P15").Select()
